$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.010.43'
$ws.Range("E2").Value = '  +0.63%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.591.42'
$ws.Range("E3").Value = '  +0.59%  '

# Row 4
$ws.Range("E4").Value = '  -0.28%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.86'
$ws.Range("E5").Value = '  +0.60%  '

# Row 6
$ws.Range("E6").Value = '  -0.26%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.480'
$ws.Range("E7").Value = '  +0.77%  '

# Row 8
$ws.Range("E8").Value = '  +0.28%  '

# Row 9
$ws.Range("E9").Value = '  -0.52%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '17.97'
$ws.Range("E10").Value = '  -0.30%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0808'
$ws.Range("E11").Value = '  +2.32%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.813.68'
$ws.Range("E12").Value = '  +0.57%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.586.73'
$ws.Range("E13").Value = '  +0.17%  '

# Row 14
$ws.Range("E14").Value = '  -0.81%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.511'
$ws.Range("E15").Value = '  +0.45%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.016.61'
$ws.Range("E16").Value = '  +0.75%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.18'
$ws.Range("E17").Value = '  +0.89%  '

# Row 18
$ws.Range("E18").Value = '  -0.23%  '

# Row 19
$ws.Range("E19").Value = '  -0.24%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '201.77'
$ws.Range("E20").Value = '  +5.28%  '

# Row 21
$ws.Range("E21").Value = '  +1.64%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.20'
$ws.Range("E22").Value = '  -1.24%  '

# Row 23
$ws.Range("E23").Value = '  +0.97%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.96'
$ws.Range("E24").Value = '  +15.74%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.19'
$ws.Range("E25").Value = '  +1.06%  '

# Row 26
$ws.Range("E26").Value = '  -0.27%  '

# Row 27
$ws.Range("E27").Value = '  -7.71%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.09'
$ws.Range("E28").Value = '  +0.15%  '

# Row 29
$ws.Range("E29").Value = '  +0.49%  '

# Row 30
$ws.Range("E30").Value = '  +0.58%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0475'
$ws.Range("E31").Value = '  +1.26%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.11'
$ws.Range("E32").Value = '  +0.10%  '

# Row 33
$ws.Range("E33").Value = '  -2.78%  '

# Row 34
$ws.Range("E34").Value = '  -0.88%  '

# Row 35
$ws.Range("E35").Value = '  -0.80%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.127.90'
$ws.Range("E36").Value = '  +2.84%  '

# Row 37
$ws.Range("E37").Value = '  +8.46%  '

# Row 38
$ws.Range("E38").Value = '  -0.21%  '

# Row 39
$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.788'
$ws.Range("E39").Value = '  +1.97%  '

# Row 40
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.31'
$ws.Range("E40").Value = '  -1.41%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.490'
$ws.Range("E41").Value = '  -1.89%  '

# Row 42
$ws.Range("E42").Value = '  -3.98%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.13'
$ws.Range("E43").Value = '  -0.43%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.724.37'
$ws.Range("E44").Value = '  +0.34%  '

# Row 45
$ws.Range("E45").Value = '  -1.66%  '

# Row 46
$ws.Range("E46").Value = '  -0.54%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '53.51'
$ws.Range("E47").Value = '  +0.87%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0503'
$ws.Range("E48").Value = '  -1.07%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.406'
$ws.Range("E49").Value = '  -0.40%  '

# Row 50
$ws.Range("E50").Value = '  +0.08%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₇0918'
$ws.Range("E51").Value = '  -17.29%  '
